$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.946.69"
$ws.Cells.Item(3, 4).Value = "1.843.11"
$ws.Cells.Item(3, 5).Value = "  +1.88%  "
$ws.Cells.Item(4, 5).Value = "  +0.15%  "
$ws.Cells.Item(5, 4).Value = "232.32"
$ws.Cells.Item(5, 5).Value = "  -0.22%  "
$ws.Cells.Item(6, 4).Value = "0.619"
$ws.Cells.Item(6, 5).Value = "  +1.22%  "
$ws.Cells.Item(7, 5).Value = "  +0.16%  "
$ws.Cells.Item(8, 4).Value = "39.83"
$ws.Cells.Item(8, 5).Value = "  -1.59%  "
$ws.Cells.Item(9, 4).Value = "0.331"
$ws.Cells.Item(9, 5).Value = "  +1.67%  "
$ws.Cells.Item(10, 4).Value = "0.0687"
$ws.Cells.Item(10, 5).Value = "  +0.52%  "
$ws.Cells.Item(11, 4).Value = "0.0983"
$ws.Cells.Item(11, 5).Value = "  -1.70%  "
$ws.Cells.Item(12, 4).Value = "2.109.84"
$ws.Cells.Item(12, 5).Value = "  +1.90%  "
$ws.Cells.Item(13, 4).Value = "11.57"
$ws.Cells.Item(13, 5).Value = "  +4.68%  "
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.848.76"
$ws.Cells.Item(14, 5).Value = "  +2.18%  "
$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(15, 4).Value = "0.676"
$ws.Cells.Item(15, 5).Value = "  +1.90%  "
$ws.Cells.Item(16, 2).Value = "Polkadot"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(16, 4).Value = "4.64"
$ws.Cells.Item(16, 5).Value = "  -0.49%  "
$ws.Cells.Item(17, 4).Value = "34.962.68"
$ws.Cells.Item(17, 5).Value = "  -0.18%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "69.90"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.27%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0787"
$ws.Cells.Item(19, 5).Value = "  -0.05%  "
$ws.Cells.Item(20, 4).Value = "239.83"
$ws.Cells.Item(21, 4).Value = "12.18"
$ws.Cells.Item(21, 5).Value = "  +1.88%  "
$ws.Cells.Item(22, 4).Value = "4.68"
$ws.Cells.Item(22, 5).Value = "  -0.79%  "
$ws.Cells.Item(23, 5).Value = "  +0.03%  "
$ws.Cells.Item(24, 4).Value = "2.29"
$ws.Cells.Item(24, 5).Value = "  +1.88%  "
$ws.Cells.Item(25, 4).Value = "171.91"
$ws.Cells.Item(25, 5).Value = "  -0.06%  "
$ws.Cells.Item(26, 4).Value = "7.81"
$ws.Cells.Item(26, 5).Value = "  -1.15%  "
$ws.Cells.Item(27, 4).Value = "17.49"
$ws.Cells.Item(27, 5).Value = "  -0.17%  "
$ws.Cells.Item(28, 5).Value = "  +2.14%  "
$ws.Cells.Item(29, 4).Value = "1.53"
$ws.Cells.Item(29, 5).Value = "  -2.37%  "
$ws.Cells.Item(30, 5).Value = "  +0.19%  "
$ws.Cells.Item(31, 4).Value = "0.0552"
$ws.Cells.Item(31, 5).Value = "  -1.06%  "
$ws.Cells.Item(32, 4).Value = "3.95"
$ws.Cells.Item(32, 5).Value = "  -3.99%  "
$ws.Cells.Item(33, 4).Value = "3.96"
$ws.Cells.Item(33, 5).Value = "  -1.76%  "
$ws.Cells.Item(34, 5).Value = "  +8.80%  "
$ws.Cells.Item(35, 4).Value = "1.23"
$ws.Cells.Item(35, 5).Value = "  +7.57%  "
$ws.Cells.Item(36, 4).Value = "1.48"
$ws.Cells.Item(36, 5).Value = "  +15.48%  "
$ws.Cells.Item(37, 4).Value = "0.702"
$ws.Cells.Item(37, 5).Value = "  +0.13%  "
$ws.Cells.Item(38, 5).Value = "  +7.30%  "
$ws.Cells.Item(39, 4).Value = "90.41"
$ws.Cells.Item(39, 5).Value = "  -2.10%  "
$ws.Cells.Item(40, 4).Value = "1.348.92"
$ws.Cells.Item(40, 5).Value = "  +2.51%  "
$ws.Cells.Item(41, 4).Value = "0.0195"
$ws.Cells.Item(41, 5).Value = "  +0.62%  "
$ws.Cells.Item(42, 4).Value = "14.84"
$ws.Cells.Item(42, 5).Value = "  +2.71%  "
$ws.Cells.Item(43, 4).Value = "2.29"
$ws.Cells.Item(43, 5).Value = "  +1.72%  "
$ws.Cells.Item(44, 5).Value = "  -2.21%  "
$ws.Cells.Item(45, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 4).Value = "0.0524"
$ws.Cells.Item(46, 5).Value = "  +2.44%  "
$ws.Cells.Item(47, 5).Value = "  -0.98%  "
$ws.Cells.Item(48, 4).Value = "2.025.25"
$ws.Cells.Item(48, 5).Value = "  +1.90%  "
$ws.Cells.Item(49, 4).Value = "3.42"
$ws.Cells.Item(49, 5).Value = "  +21.25%  "
$ws.Cells.Item(50, 5).Value = "  +0.18%  "
$ws.Cells.Item(51, 4).Value = "0.0668"
$ws.Cells.Item(51, 5).Value = "  +0.08%  "
